$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.317.97"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.843.72"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'233.33"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.4664"
$ws.Range("E7").Value = "  -1.82%  "
$ws.Range("D8").Value = "'0.2732"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.06287"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").Value = "1.838.87"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'0.07442"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "'16.28"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "'4.943"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "'83.84"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "'0.6207"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").Value = "30.270.93"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'0.9993"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'228.11"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").Value = "'12.38"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'4.917"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").Value = "'5.866"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").Value = "'9.188"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'164.70"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "'17.82"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'1.875"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'0.1029"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'1.372"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'4.076"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'3.806"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("D32").Value = "'0.04852"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'1.143"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "'0.7118"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "'2.701"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "'0.01899"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'2.654"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "'0.8850"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").Value = "'105.50"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'1.924"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'5.555"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'0.4028"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("D44").Value = "'7.139"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").Value = "'62.22"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").Value = "'0.1199"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "'8.600"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "'33.23"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "'0.05510"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'1.353"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").Value = "'0.3641"
$ws.Range("E51").Value = "  -1.65%  "
